$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 15 de Octubre de 2020 a las 15:26"

# Reorder countries: Eslovenia moves up (rows 107-112 shift), data refreshed
$ws.Range("A107").Value = "Eslovenia"
$ws.Range("A108").Value = "Mozambique"
$ws.Range("A109").Value = "Tayikistan"
$ws.Range("A110").Value = "Guayana Francesa"
$ws.Range("A111").Value = "Uganda"
$ws.Range("A112").Value = "Luxemburgo"

# Update numeric statistics that changed
$ws.Cells.Item(4, 2).Value = 8156124
$ws.Cells.Item(4, 3).Value = 6081
$ws.Cells.Item(4, 4).Value = 5280004
$ws.Cells.Item(4, 5).Value = 2654225
$ws.Cells.Item(4, 7).Value = 52
$ws.Cells.Item(4, 8).Value = 221895

$ws.Cells.Item(24, 2).Value = 341062
$ws.Cells.Item(24, 3).Value = 472
$ws.Cells.Item(24, 4).Value = 327327
$ws.Cells.Item(24, 5).Value = 8608
$ws.Cells.Item(24, 7).Value = 19
$ws.Cells.Item(24, 8).Value = 5127

$ws.Cells.Item(29, 2).Value = 203954
$ws.Cells.Item(29, 3).Value = 7791
$ws.Cells.Item(29, 7).Value = 29
$ws.Cells.Item(29, 8).Value = 6692

$ws.Cells.Item(38, 2).Value = 128803
$ws.Cells.Item(38, 3).Value = 200
$ws.Cells.Item(38, 4).Value = 125802
$ws.Cells.Item(38, 5).Value = 2779
$ws.Cells.Item(38, 7).Value = 2
$ws.Cells.Item(38, 8).Value = 222

$ws.Cells.Item(43, 2).Value = 111437
$ws.Cells.Item(43, 3).Value = 1398
$ws.Cells.Item(43, 4).Value = 103325
$ws.Cells.Item(43, 5).Value = 7660
$ws.Cells.Item(43, 7).Value = 2
$ws.Cells.Item(43, 8).Value = 452

$ws.Cells.Item(47, 2).Value = 102407
$ws.Cells.Item(47, 7).Value = 6
$ws.Cells.Item(47, 8).Value = 5910

$ws.Cells.Item(57, 5).Value = 3772
$ws.Cells.Item(57, 7).Value = 1
$ws.Cells.Item(57, 8).Value = 288

$ws.Cells.Item(58, 4).Value = 50500
$ws.Cells.Item(58, 5).Value = 18704
$ws.Cells.Item(58, 7).Value = 4
$ws.Cells.Item(58, 8).Value = 2113

$ws.Cells.Item(70, 2).Value = 46676
$ws.Cells.Item(70, 3).Value = 855
$ws.Cells.Item(70, 4).Value = 25685
$ws.Cells.Item(70, 5).Value = 20310
$ws.Cells.Item(70, 7).Value = 12
$ws.Cells.Item(70, 8).Value = 681

$ws.Cells.Item(77, 2).Value = 35454
$ws.Cells.Item(77, 3).Value = 203
$ws.Cells.Item(77, 5).Value = 3148
$ws.Cells.Item(77, 7).Value = 2
$ws.Cells.Item(77, 8).Value = 770

$ws.Cells.Item(80, 2).Value = 32224
$ws.Cells.Item(80, 3).Value = 569
$ws.Cells.Item(80, 4).Value = 24356
$ws.Cells.Item(80, 5).Value = 6896
$ws.Cells.Item(80, 7).Value = 14
$ws.Cells.Item(80, 8).Value = 972

$ws.Cells.Item(96, 5).Value = 3812
$ws.Cells.Item(96, 7).Value = 1
$ws.Cells.Item(96, 8).Value = 278

$ws.Cells.Item(107, 2).Value = 10683
$ws.Cells.Item(107, 3).Value = 745
$ws.Cells.Item(107, 4).Value = 5689
$ws.Cells.Item(107, 5).Value = 4818
$ws.Cells.Item(107, 7).Value = 1
$ws.Cells.Item(107, 8).Value = 176

$ws.Cells.Item(108, 2).Value = 10392
$ws.Cells.Item(108, 4).Value = 8035
$ws.Cells.Item(108, 5).Value = 2284
$ws.Cells.Item(108, 8).Value = 73

$ws.Cells.Item(109, 2).Value = 10374
$ws.Cells.Item(109, 3).Value = 38
$ws.Cells.Item(109, 4).Value = 9317
$ws.Cells.Item(109, 5).Value = 977
$ws.Cells.Item(109, 7).Value = 1
$ws.Cells.Item(109, 8).Value = 80

$ws.Cells.Item(110, 2).Value = 10202
$ws.Cells.Item(110, 4).Value = 9892
$ws.Cells.Item(110, 5).Value = 241
$ws.Cells.Item(110, 8).Value = 69

$ws.Cells.Item(111, 2).Value = 10069
$ws.Cells.Item(111, 4).Value = 6531
$ws.Cells.Item(111, 5).Value = 3443
$ws.Cells.Item(111, 8).Value = 95

$ws.Cells.Item(112, 2).Value = 10030
$ws.Cells.Item(112, 4).Value = 8306
$ws.Cells.Item(112, 5).Value = 1591
$ws.Cells.Item(112, 8).Value = 133

$ws.Cells.Item(128, 2).Value = 5214
$ws.Cells.Item(128, 3).Value = 12
$ws.Cells.Item(128, 4).Value = 4943
$ws.Cells.Item(128, 5).Value = 166

